$d = $word.ActiveDocument

# The document currently contains a single paragraph "test" that carries
# the hidden "_GoBack" bookmark right after the word. The target revision
# prepends a new paragraph "DON’T FORGET TO COMMIT AND PUSH!!!", followed
# by a blank paragraph, and moves the "_GoBack" bookmark so that it now
# sits before the "test" run instead of after it.

$firstPara = $d.Paragraphs(1).Range
$insertionPoint = $firstPara.Duplicate
$insertionPoint.Collapse(1)

$newContentXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>DON' + [char]0x2019 + 'T FORGET TO COMMIT AND PUSH!!!</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newContentXml) | Out-Null

# Re-find the paragraph that now holds "test" (it kept its own identity,
# so it is the third paragraph after the two we just inserted) and move
# the "_GoBack" bookmark to its very start.
$testParaRange = $d.Paragraphs(3).Range
$bookmarkTarget = $testParaRange.Duplicate
$bookmarkTarget.Collapse(1)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Add("_GoBack", $bookmarkTarget) | Out-Null
}
